$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FlagReason")

# --- New header cell C1 ("Comment 2"), formatted like the other header cells ---
$ws.Range("C1").Value = "Comment 2"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- New data cell D2 (set before D1 / B2 so shared-string indices line up with source order) ---
$ws.Range("D2").Value = "Requesting to change Company Type to Operating Company with Ownership: Private Equity Group because it is being considered to be a potential round trip"
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# --- Update existing comment text in B2 ---
$ws.Range("B2").Value = "Requesting to change Company Type to Operating Company because it is being considered to be a potential round trip"

# --- New header cell D1 ("Comment 3") ---
$ws.Range("D1").Value = "Comment 3"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# --- New data cell C2 (reuses an already-existing shared string) ---
$ws.Range("C2").Value = "Requesting to change Company Type to Operating Company and Ownership to Private Equity Group because it is being considered to be a potential round trip"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths for the two new columns (C matches existing B width, D is narrower)
$ws.Columns.Item(3).ColumnWidth = 55.17
$ws.Columns.Item(4).ColumnWidth = 46.5

# Row height grows because of the additional wrapped text
$ws.Rows.Item(2).RowHeight = 57.6

# Make FlagReason the active/selected sheet with D1 selected
$ws.Activate() | Out-Null
$ws.Range("D1").Select() | Out-Null
